$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data for rows 2-17 (A:B:C). Row 18 (Christian Braun / Denver Nuggets) is removed.
$data = @(
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Mike Conley", "PG", "Minnesota Timberwolves"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Brandon Boston Jr.", "SG,SF,PF", "New Orleans Pelicans"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused last row (row 18), which previously held Christian Braun.
$ws.Rows.Item(18).Delete()
